# harmonised offset data (04.02.2026 version)
# Rebuilds the BasePointers / GameInfo tables: drops the bold/bordered
# header style and re-lays out the rows (new pointer/field entries
# inserted alphabetically, values refreshed for the 2k26 columns).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BasePointers")
$ws2 = $wb.Worksheets.Item("GameInfo")

# ---------------------------------------------------------------------
# Drop the bold-centered/bordered header style from row 1 of both
# sheets -- headers become plain (default) styled cells.
# ---------------------------------------------------------------------
$ws1.Range("A1:F1").ClearFormats()
$ws2.Range("A1:F1").ClearFormats()

# ---------------------------------------------------------------------
# BasePointers: re-lay the pointer rows out alphabetically, with the
# new HallOfFame / History / Jersey / NBAHistory / Stadium / Staff /
# TeamHistory / career_stats pointers added and the 2k26 values
# refreshed for Player / Team.
# ---------------------------------------------------------------------
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = 'HallOfFame'
$ws1.Range("A3").ClearFormats()
$ws1.Range("E3").Value = ""
$ws1.Range("F3").NumberFormat = "@"
$ws1.Range("F3").Value = '130991424'
$ws1.Range("F3").ClearFormats()

$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = 'History'
$ws1.Range("A4").ClearFormats()
$ws1.Range("E4").Value = ""
$ws1.Range("F4").NumberFormat = "@"
$ws1.Range("F4").Value = '130991424'
$ws1.Range("F4").ClearFormats()

$ws1.Range("A5").NumberFormat = "@"
$ws1.Range("A5").Value = 'Jersey'
$ws1.Range("A5").ClearFormats()
$ws1.Range("F5").NumberFormat = "@"
$ws1.Range("F5").Value = '130991400'
$ws1.Range("F5").ClearFormats()

$ws1.Range("A6").NumberFormat = "@"
$ws1.Range("A6").Value = 'NBAHistory'
$ws1.Range("A6").ClearFormats()
$ws1.Range("F6").NumberFormat = "@"
$ws1.Range("F6").Value = '130991424'
$ws1.Range("F6").ClearFormats()

$ws1.Range("A7").NumberFormat = "@"
$ws1.Range("A7").Value = 'Player'
$ws1.Range("A7").ClearFormats()
$ws1.Range("E7").NumberFormat = "@"
$ws1.Range("E7").Value = '132504968'
$ws1.Range("E7").ClearFormats()
$ws1.Range("F7").NumberFormat = "@"
$ws1.Range("F7").Value = '130990776'
$ws1.Range("F7").ClearFormats()

$ws1.Range("A8").NumberFormat = "@"
$ws1.Range("A8").Value = 'Stadium'
$ws1.Range("A8").ClearFormats()
$ws1.Range("F8").NumberFormat = "@"
$ws1.Range("F8").Value = '130991160'
$ws1.Range("F8").ClearFormats()

$ws1.Range("A9").NumberFormat = "@"
$ws1.Range("A9").Value = 'Staff'
$ws1.Range("A9").ClearFormats()
$ws1.Range("F9").NumberFormat = "@"
$ws1.Range("F9").Value = '130991496'
$ws1.Range("F9").ClearFormats()

$ws1.Range("A10").NumberFormat = "@"
$ws1.Range("A10").Value = 'Team'
$ws1.Range("A10").ClearFormats()
$ws1.Range("E10").NumberFormat = "@"
$ws1.Range("E10").Value = '132505568'
$ws1.Range("E10").ClearFormats()
$ws1.Range("F10").NumberFormat = "@"
$ws1.Range("F10").Value = '130991376'
$ws1.Range("F10").ClearFormats()

$ws1.Range("A11").NumberFormat = "@"
$ws1.Range("A11").Value = 'TeamHistory'
$ws1.Range("A11").ClearFormats()
$ws1.Range("F11").NumberFormat = "@"
$ws1.Range("F11").Value = '130991376'
$ws1.Range("F11").ClearFormats()

$ws1.Range("A12").NumberFormat = "@"
$ws1.Range("A12").Value = 'career_stats'
$ws1.Range("A12").ClearFormats()
$ws1.Range("F12").NumberFormat = "@"
$ws1.Range("F12").Value = '130990680'
$ws1.Range("F12").ClearFormats()

# ---------------------------------------------------------------------
# GameInfo: re-lay the field rows out alphabetically, with the new
# career_statsSize / coachSize / hall_of_fameSize / historySize /
# jerseySize / stadiumSize / staffSize fields added, and refresh the
# executable / size / version columns for 2k25 & 2k26.
# ---------------------------------------------------------------------
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = 'career_statsSize'
$ws2.Range("A2").ClearFormats()
$ws2.Range("B2").Value = ""
$ws2.Range("C2").Value = ""
$ws2.Range("D2").Value = ""
$ws2.Range("E2").Value = ""
$ws2.Range("F2").Value = 64

$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = 'coachSize'
$ws2.Range("A3").ClearFormats()
$ws2.Range("E3").Value = ""
$ws2.Range("F3").Value = 432

$ws2.Range("A4").NumberFormat = "@"
$ws2.Range("A4").Value = 'executable'
$ws2.Range("A4").ClearFormats()
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = 'NBA2K22.exe'
$ws2.Range("B4").ClearFormats()
$ws2.Range("C4").NumberFormat = "@"
$ws2.Range("C4").Value = 'NBA2K23.exe'
$ws2.Range("C4").ClearFormats()
$ws2.Range("D4").NumberFormat = "@"
$ws2.Range("D4").Value = 'NBA2K24.exe'
$ws2.Range("D4").ClearFormats()
$ws2.Range("E4").NumberFormat = "@"
$ws2.Range("E4").Value = 'NBA2K25.exe'
$ws2.Range("E4").ClearFormats()
$ws2.Range("F4").NumberFormat = "@"
$ws2.Range("F4").Value = 'NBA2K26.exe'
$ws2.Range("F4").ClearFormats()

$ws2.Range("A5").NumberFormat = "@"
$ws2.Range("A5").Value = 'hall_of_fameSize'
$ws2.Range("A5").ClearFormats()
$ws2.Range("E5").Value = ""
$ws2.Range("F5").Value = 108

$ws2.Range("A6").NumberFormat = "@"
$ws2.Range("A6").Value = 'historySize'
$ws2.Range("A6").ClearFormats()
$ws2.Range("F6").Value = 168

$ws2.Range("A7").NumberFormat = "@"
$ws2.Range("A7").Value = 'jerseySize'
$ws2.Range("A7").ClearFormats()
$ws2.Range("F7").Value = 368

$ws2.Range("A8").NumberFormat = "@"
$ws2.Range("A8").Value = 'playerSize'
$ws2.Range("A8").ClearFormats()
$ws2.Range("E8").Value = 1096
$ws2.Range("F8").Value = 1176

$ws2.Range("A9").NumberFormat = "@"
$ws2.Range("A9").Value = 'stadiumSize'
$ws2.Range("A9").ClearFormats()
$ws2.Range("F9").Value = 4792

$ws2.Range("A10").NumberFormat = "@"
$ws2.Range("A10").Value = 'staffSize'
$ws2.Range("A10").ClearFormats()
$ws2.Range("F10").Value = 432

$ws2.Range("A11").NumberFormat = "@"
$ws2.Range("A11").Value = 'teamSize'
$ws2.Range("A11").ClearFormats()
$ws2.Range("E11").Value = 5664
$ws2.Range("F11").Value = 5672

$ws2.Range("A12").NumberFormat = "@"
$ws2.Range("A12").Value = 'version'
$ws2.Range("A12").ClearFormats()
$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = '2K22'
$ws2.Range("B12").ClearFormats()
$ws2.Range("C12").NumberFormat = "@"
$ws2.Range("C12").Value = '2K23'
$ws2.Range("C12").ClearFormats()
$ws2.Range("D12").NumberFormat = "@"
$ws2.Range("D12").Value = '2K24'
$ws2.Range("D12").ClearFormats()
$ws2.Range("E12").NumberFormat = "@"
$ws2.Range("E12").Value = '2K25'
$ws2.Range("E12").ClearFormats()
$ws2.Range("F12").NumberFormat = "@"
$ws2.Range("F12").Value = '2026-02-21'
$ws2.Range("F12").ClearFormats()
